# Add "Quantity" (column F) and "Price" (column G) columns to the PRODUCTS sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PRODUCTS")

# Header cells - values first, then copy the header style/format from an
# existing header cell (D1) so F1/G1 match the look of the other headers.
$ws.Range("F1").Value = "Quantity"
$ws.Range("G1").Value = "Price"
$ws.Range("D1").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Data values for rows 2..31 (Quantity / Price), in sheet order.
$fVals = @(0,0,0,0,0,0,0,0,0,0,0,0,2,0,0,0,1,0,0,0,0,0,0,1,0,0,0,0,0,0)
$gVals = @(1999,3299,1799,999,1999,449,499.99,799.99,0,0,0,25.99,15.99,13.99,23.99,23.99,23.99,23.99,33.450000000000003,23.99,48.59,23.99,237,523,229.95,49.95,599,149,828,9.99)

for ($i = 0; $i -lt $fVals.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $fVals[$i]
    $ws.Cells.Item($row, 7).Value = $gVals[$i]
}

# Match the selection shown in the saved workbook.
$ws.Range("G1").Select() | Out-Null
